$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update the query timestamps in the "data" sheet (F2:F27) ---
$data.Range("F2").Value  = "2021-10-05 14:21:18.537095"
$data.Range("F3").Value  = "2021-10-05 14:21:18.537101"
$data.Range("F4").Value  = "2021-10-05 14:21:18.537104"
$data.Range("F5").Value  = "2021-10-05 14:21:18.537106"
$data.Range("F6").Value  = "2021-10-05 14:21:18.537108"
$data.Range("F7").Value  = "2021-10-05 14:21:18.537110"
$data.Range("F8").Value  = "2021-10-05 14:21:18.537112"
$data.Range("F9").Value  = "2021-10-05 14:21:18.537114"
$data.Range("F10").Value = "2021-10-05 14:21:18.537116"
$data.Range("F11").Value = "2021-10-05 14:21:18.537119"
$data.Range("F12").Value = "2021-10-05 14:21:18.537122"
$data.Range("F13").Value = "2021-10-05 14:21:18.537125"
$data.Range("F14").Value = "2021-10-05 14:21:18.537128"
$data.Range("F15").Value = "2021-10-05 14:21:18.537131"
$data.Range("F16").Value = "2021-10-05 14:21:18.537133"
$data.Range("F17").Value = "2021-10-05 14:21:18.537136"
$data.Range("F18").Value = "2021-10-05 14:21:18.537138"
$data.Range("F19").Value = "2021-10-05 14:21:18.537140"
$data.Range("F20").Value = "2021-10-05 14:21:18.537142"
$data.Range("F21").Value = "2021-10-05 14:21:18.537144"
$data.Range("F22").Value = "2021-10-05 14:21:18.537146"
$data.Range("F23").Value = "2021-10-05 14:21:18.537148"
$data.Range("F24").Value = "2021-10-05 14:21:18.537149"
$data.Range("F25").Value = "2021-10-05 14:21:18.537151"
$data.Range("F26").Value = "2021-10-05 14:21:18.537180"
$data.Range("F27").Value = "2021-10-05 14:21:18.537182"

# --- Add the new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Reuse the header style (bold + border + centered) from the "data" sheet header row
$data.Range("B1:F1").Copy($meta.Range("B1:F1"))
$data.Range("F1").Copy($meta.Range("G1"))
# Reuse the numeric-id style from the "data" sheet
$data.Range("A2").Copy($meta.Range("A2"))

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Iron metabolism disorders"
$meta.Range("C2").Value = 515
# Keep the version number as text ("1.33"), not a float
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.33"
$meta.Range("E2").Value = "2021-03-03T10:27:31.439328Z"
$meta.Range("F2").Value = "2021-10-05 14:21:18.534637"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/515/?format=json"

# Match the page margins used elsewhere in the workbook (inches -> points: *72)
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36

# Keep "data" as the active/selected sheet (matches activeTab=0 in the original)
$data.Activate()
